$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "286.71"
Set-TextValue "E2" "4.16%"
Set-TextValue "D3" "28.35"
Set-TextValue "E3" "4.26%"
Set-TextValue "D4" "4.934"
Set-TextValue "E4" "1.96%"
Set-TextValue "D5" "0.06550"
Set-TextValue "E5" "2.45%"
Set-TextValue "D6" "7.248"
Set-TextValue "E6" "4.39%"
Set-TextValue "D7" "1.342"
Set-TextValue "E7" "10.14%"
Set-TextValue "D8" "0.9181"
Set-TextValue "E8" "4.80%"
Set-TextValue "D9" "0.1572"
Set-TextValue "E9" "3.57%"
Set-TextValue "D10" "0.06549"
Set-TextValue "E10" "29.10%"
Set-TextValue "D11" "0.07653"
Set-TextValue "E11" "1.79%"
Set-TextValue "D12" "0.02981"
Set-TextValue "E12" "0.83%"
Set-TextValue "D13" "0.08973"
Set-TextValue "E13" "-0.25%"
Set-TextValue "E14" "2.15%"
Set-TextValue "D15" "0.0006545"
Set-TextValue "E15" "2.10%"
Set-TextValue "D16" "0.006071"
Set-TextValue "E16" "-1.80%"
Set-TextValue "D17" "3.488"
Set-TextValue "E17" "0.66%"
Set-TextValue "D18" "3.385"
Set-TextValue "E18" "2.23%"
Set-TextValue "E19" "-1.34%"
Set-TextValue "D21" "0.1349"
Set-TextValue "E21" "0.58%"
Set-TextValue "D22" "4.003"
Set-TextValue "E22" "2.54%"
Set-TextValue "D24" "0.04467"
Set-TextValue "E24" "0.96%"
Set-TextValue "E25" "0.90%"
Set-TextValue "D26" "0.004343"
Set-TextValue "E26" "12.75%"
Set-TextValue "E28" "-1.80%"
Set-TextValue "E29" "-15.79%"
Set-TextValue "D40" "0.04164"
Set-TextValue "E40" "0.89%"
Set-TextValue "D41" "0.006886"
Set-TextValue "E41" "1.29%"
Set-TextValue "E42" "20.27%"
Set-TextValue "D43" "0.002058"
Set-TextValue "E43" "-3.87%"
Set-TextValue "D44" "0.01249"
Set-TextValue "E44" "8.90%"
Set-TextValue "D45" "0.00005548"
Set-TextValue "E45" "7.19%"
Set-TextValue "D46" "1.561"
Set-TextValue "E46" "3.97%"
Set-TextValue "D47" "0.01849"
Set-TextValue "E47" "-7.63%"
